# Generate Report for Handback
# Adds a new handback row (cedf3195-1b0e-4132-aff5-976f0a3e29ef) to the
# Overview / zh-cn / de-de sheets of the handback-status report.

$wb = $excel.ActiveWorkbook

$fileId   = "cedf3195-1b0e-4132-aff5-976f0a3e29ef"
$fileHash = "c9e94b16ece26f7a6e72a424fc9b03188a77d09b"
$mdName   = "$fileId.md"
$zhName   = "$fileId.$fileHash.zh-cn.xlf"
$deName   = "$fileId.$fileHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$includeText  = "Include"

# VBA/Excel xlUnderlineStyleSingle = 2
$underlineSingle = 2
# RGB(0x64,0x95,0xED) packed as BGR (matches the workbook's existing HyperLink font color FF6495ED)
$linkColor = 15570276

function Style-AsLink($range) {
    $range.Font.Underline = $underlineSingle
    $range.Font.Color = $linkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
Style-AsLink $wsOverview.Range("A4")

$wsOverview.Range("B4").Value2 = $statusInSync
$wsOverview.Range("C4").Value2 = $statusInSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
Style-AsLink $wsZh.Range("A4")

$wsZh.Range("B4").Value2 = $statusInSync

$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$fileId/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhName",
    "",
    "",
    $zhName) | Out-Null
Style-AsLink $wsZh.Range("C4")

$wsZh.Range("D4").Value2 = "2016-02-26 06:20:05"
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$fileId/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
Style-AsLink $wsZh.Range("E4")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$fileId/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhName",
    "",
    "",
    $zhName) | Out-Null
Style-AsLink $wsZh.Range("F4")

$wsZh.Range("G4").Value2 = "2016-02-26 06:20:51"
$wsZh.Range("H4").Value2 = $includeText

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
Style-AsLink $wsDe.Range("A4")

$wsDe.Range("B4").Value2 = $statusInSync

$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$fileId/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deName",
    "",
    "",
    $deName) | Out-Null
Style-AsLink $wsDe.Range("C4")

$wsDe.Range("D4").Value2 = "2016-02-26 06:20:18"
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$fileId/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
Style-AsLink $wsDe.Range("E4")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$fileId/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deName",
    "",
    "",
    $deName) | Out-Null
Style-AsLink $wsDe.Range("F4")

$wsDe.Range("G4").Value2 = "2016-02-26 06:21:14"
$wsDe.Range("H4").Value2 = $includeText

Write-Output "Handback report row added for $fileId"
